$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts per diff (column B = Encontrados_GitHub, column C = Encontrados_GitLab)
$ws.Range("B2").Value = 2
# C2 unchanged (0)

$ws.Range("B3").Value = 39
$ws.Range("C3").Value = 0

$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 0

$ws.Range("B5").Value = 148
# C5 unchanged (0)

$ws.Range("B6").Value = 5
# C6 unchanged (0)

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0

$ws.Range("B11").Value = 1
# C11 unchanged (0)

$ws.Range("B12").Value = 3
# C12 unchanged (0)

$ws.Range("B15").Value = 176
$ws.Range("C15").Value = 0
